$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.31185386054396
$ws.Range("C2").Value = 10.11626123327131
$ws.Range("E2").Value = 12.02127045850227
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 15.85889875550464
$ws.Range("H2").Value = 11.13038194287012
$ws.Range("M2").Value = 14.08220036599496
$ws.Range("O2").Value = 15.1985073221303
$ws.Range("B3").Value = 11.59983846277207
$ws.Range("C3").Value = 9.748213023181888
$ws.Range("E3").Value = 11.97245278535159
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 16.02954982958105
$ws.Range("H3").Value = 11.2022018525049
$ws.Range("M3").Value = 13.7322606534324
$ws.Range("O3").Value = 15.32945736226491
$ws.Range("B4").Value = 11.13874100402733
$ws.Range("C4").Value = 9.514618661399259
$ws.Range("E4").Value = 11.94829259168963
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 16.1464153507662
$ws.Range("H4").Value = 11.2489548823832
$ws.Range("M4").Value = 13.51488933768148
$ws.Range("O4").Value = 15.41542754251826
$ws.Range("B5").Value = 10.94492874164935
$ws.Range("C5").Value = 9.417631514144853
$ws.Range("E5").Value = 11.93991396945585
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 16.1970330931781
$ws.Range("H5").Value = 11.26867483665019
$ws.Range("M5").Value = 13.42580576585454
$ws.Range("O5").Value = 15.45185578075706
$ws.Range("B6").Value = 10.9123928214129
$ws.Range("C6").Value = 9.401422140667803
$ws.Range("E6").Value = 11.93861138478871
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 16.20561765805988
$ws.Range("H6").Value = 11.27198964964152
$ws.Range("M6").Value = 13.41098696476583
$ws.Range("O6").Value = 15.45798875492476
$ws.Range("B7").Value = 11.13615096516846
$ws.Range("C7").Value = 9.513317762949214
$ws.Range("E7").Value = 11.94817365118775
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 16.14708593556965
$ws.Range("H7").Value = 11.24921812899067
$ws.Range("M7").Value = 13.51368978944069
$ws.Range("O7").Value = 15.41591318593436
$ws.Range("B8").Value = 12.07136952402198
$ws.Range("C8").Value = 9.991004962637192
$ws.Range("E8").Value = 12.0032352262813
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 15.91520405406711
$ws.Range("H8").Value = 11.15459433924502
$ws.Range("M8").Value = 13.96213600829622
$ws.Range("O8").Value = 15.24250054988653
$ws.Range("B9").Value = 13.75421859045071
$ws.Range("C9").Value = 10.86279749348028
$ws.Range("E9").Value = 12.15698725062742
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 15.55836978930097
$ws.Range("H9").Value = 10.9901056053752
$ws.Range("M9").Value = 14.81637159306386
$ws.Range("O9").Value = 14.94683537137701
$ws.Range("B10").Value = 14.87794316594925
$ws.Range("C10").Value = 11.45843910347277
$ws.Range("E10").Value = 12.29717968704845
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 15.3585809488787
$ws.Range("H10").Value = 10.88209921434705
$ws.Range("M10").Value = 15.42206788461447
$ws.Range("O10").Value = 14.75700079558623
$ws.Range("B11").Value = 15.35995621029956
$ws.Range("C11").Value = 11.71877368985363
$ws.Range("E11").Value = 12.36666833736408
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 15.28180312738522
$ws.Range("H11").Value = 10.83575414009229
$ws.Range("M11").Value = 15.69161766836012
$ws.Range("O11").Value = 14.67666205120963
$ws.Range("B12").Value = 15.53828789678284
$ws.Range("C12").Value = 11.81576562966889
$ws.Range("E12").Value = 12.39378293855381
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 15.25480322122782
$ws.Range("H12").Value = 10.81860560717687
$ws.Range("M12").Value = 15.792735761312
$ws.Range("O12").Value = 14.6471118734113
$ws.Range("B13").Value = 15.50006757818447
$ws.Range("C13").Value = 11.79494836070799
$ws.Range("E13").Value = 12.38790806818154
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 15.26052517251691
$ws.Range("H13").Value = 10.82228099236212
$ws.Range("M13").Value = 15.77100203908445
$ws.Range("O13").Value = 14.65343712060126
$ws.Range("B14").Value = 15.37471175504523
$ws.Range("C14").Value = 11.72678549494094
$ws.Range("E14").Value = 12.36888315175811
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 15.27954000266484
$ws.Range("H14").Value = 10.83433527347397
$ws.Range("M14").Value = 15.69995635628486
$ws.Range("O14").Value = 14.67421341600589
$ws.Range("B15").Value = 15.29738126453512
$ws.Range("C15").Value = 11.68482483267861
$ws.Range("E15").Value = 12.35733346198957
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 15.29145858615625
$ws.Range("H15").Value = 10.84177115074136
$ws.Range("M15").Value = 15.65631185542438
$ws.Range("O15").Value = 14.68705331547889
$ws.Range("B16").Value = 14.84585524222085
$ws.Range("C16").Value = 11.44120621306033
$ws.Range("E16").Value = 12.29275168944884
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 15.36388635932314
$ws.Range("H16").Value = 10.88518410768436
$ws.Range("M16").Value = 15.40432353580259
$ws.Range("O16").Value = 14.76237280984663
$ws.Range("B17").Value = 14.5613843699156
$ws.Range("C17").Value = 11.28898486371706
$ws.Range("E17").Value = 12.25458228677711
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 15.4119644721388
$ws.Range("H17").Value = 10.91253084591324
$ws.Range("M17").Value = 15.24813473284359
$ws.Range("O17").Value = 14.81012511765762
$ws.Range("B18").Value = 14.39501794062377
$ws.Range("C18").Value = 11.20043421550608
$ws.Range("E18").Value = 12.23316753285209
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 15.44094290359211
$ws.Range("H18").Value = 10.92852231306924
$ws.Range("M18").Value = 15.15774089672676
$ws.Range("O18").Value = 14.83815698777682
$ws.Range("B19").Value = 14.33821753667796
$ws.Range("C19").Value = 11.17028332615204
$ws.Range("E19").Value = 12.22601008710644
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 15.45098069185971
$ws.Range("H19").Value = 10.93398179339431
$ws.Range("M19").Value = 15.12704222310338
$ws.Range("O19").Value = 14.84774510046467
$ws.Range("B20").Value = 14.59195103085342
$ws.Range("C20").Value = 11.30529275192422
$ws.Range("E20").Value = 12.25858981624768
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 15.4067089900539
$ws.Range("H20").Value = 10.90959258367227
$ws.Range("M20").Value = 15.2648197714271
$ws.Range("O20").Value = 14.80498317364473
$ws.Range("B21").Value = 15.41164567069918
$ws.Range("C21").Value = 11.74685022445076
$ws.Range("E21").Value = 12.37444967696485
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 15.27389823656605
$ws.Range("H21").Value = 10.83078374242332
$ws.Range("M21").Value = 15.72085081252474
$ws.Range("O21").Value = 14.66808718191296
$ws.Range("B22").Value = 15.92290698586819
$ws.Range("C22").Value = 12.02613584312096
$ws.Range("E22").Value = 12.45482693133306
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 15.19920881437082
$ws.Range("H22").Value = 10.7816172015251
$ws.Range("M22").Value = 16.01329140653004
$ws.Range("O22").Value = 14.58370593953788
$ws.Range("B23").Value = 15.65227432056921
$ws.Range("C23").Value = 11.87794529450097
$ws.Range("E23").Value = 12.41150954691901
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 15.23794919298513
$ws.Range("H23").Value = 10.80764405896523
$ws.Range("M23").Value = 15.85775199726979
$ws.Range("O23").Value = 14.62827385484875
$ws.Range("B24").Value = 14.57814062811144
$ws.Range("C24").Value = 11.29792317949248
$ws.Range("E24").Value = 12.25677636175055
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 15.40908083018208
$ws.Range("H24").Value = 10.91092013268646
$ws.Range("M24").Value = 15.25727832342706
$ws.Range("O24").Value = 14.80730604524692
$ws.Range("B25").Value = 13.31455670207033
$ws.Range("C25").Value = 10.63449542071862
$ws.Range("E25").Value = 12.11055132678083
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 15.64414274121531
$ws.Range("H25").Value = 11.03234890016742
$ws.Range("M25").Value = 14.5886794547111
$ws.Range("O25").Value = 15.02203492316789
